# Insert a new "Match ID" column at the front of the sheet (shifts A:AC -> B:AD)
# and populate it: header label in row 3, match id 21 for every data row (4-19),
# and 21 for the hidden totals row (20) without the bold header styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing column one place to the right.
$ws.Columns.Item(1).Insert() | Out-Null

# Row 3 holds the (visible) column header labels - add the new "Match ID" label,
# bold to match the other header cells' font (no border, unlike the B3:AD3 header row).
$ws.Range("A3").Value2 = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Data rows 4-19: same match id for every player, bold like the header style that
# was minted for column A.
for ($r = 4; $r -le 19; $r++) {
    $ws.Range("A$r").Value2 = 21
    $ws.Range("A$r").Font.Bold = $true
}

# Row 20 is the hidden totals row - gets the match id too, but keeps the default
# (non-bold) style, matching the rest of that row's formatting.
$ws.Range("A20").Value2 = 21

# Restore the selection to match the new layout.
$ws.Range("A3:A19").Select() | Out-Null
